$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new columns (P1, Q1) continuing the sequence 0..15
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header cell formatting (bold, centered, bordered) from an existing
# header cell (B1) onto the two new header cells without disturbing values.
$ws.Range("B1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update data rows 2-25: swap values in columns I/K and M/O, and populate new P,Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column, all 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column, all 2
}
